# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Swap the table order of a few countries whose rows leapfrogged each
#   other in the source feed (Benin/Jamaica, Dominica/Fiyi, Islas Turcas y
#   Caicos/Santa Sede), carrying each country's own data with it
# - Refresh the case/recovered/death counters for the countries whose
#   stats moved between scrapes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 04:34"

# Bolivia (row 49) - updated counters
$ws.Range("B49").Value = 22476
$ws.Range("C49").Value = 977
$ws.Range("D49").Value = 4670
$ws.Range("E49").Value = 17091
$ws.Range("G49").Value = 18
$ws.Range("H49").Value = 715

# Corea del Sur (row 62) - updated counters
$ws.Range("B62").Value = 12373
$ws.Range("C62").Value = 67
$ws.Range("D62").Value = 10856
$ws.Range("E62").Value = 1237

# Nueva Zelanda (row 118) - updated counters
$ws.Range("B118").Value = 1509
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 5

# Rows 145/146: Benin and Jamaica swap places, each keeping its own
# (refreshed) data
$ws.Range("A145").Value = "Jamaica"
$ws.Range("B145").Value = 652
$ws.Range("C145").Value = 14
$ws.Range("D145").Value = 458
$ws.Range("E145").Value = 184
$ws.Range("H145").Value = 10

$ws.Range("A146").Value = "Benin"
$ws.Range("B146").Value = 650
$ws.Range("D146").Value = 247
$ws.Range("E146").Value = 392
$ws.Range("H146").Value = 11

# Bermudas (row 172) - updated counters
$ws.Range("B172").Value = 146
$ws.Range("C172").Value = 2
$ws.Range("D172").Value = 132
$ws.Range("E172").Value = 5

# Rows 202/203: Dominica and Fiyi swap places (identical data either way)
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Rows 208/209: Islas Turcas y Caicos and Santa Sede swap places, each
# keeping its own data
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
